$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the view: scroll the window so row 7 is at the top-left, and
#     move the active selection to J36 (best-effort; the headless engine's
#     view-state model only persists the selection, not the scroll offset). ---
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

# --- Column I: base-shear (D6, D7, ...) divided by the matching D value
#     for rows 23-37 ("base shear to Sa" ratio). ---
$ws.Cells.Item(23, 9).Formula = "=D6/D23"
$iFillRange = $ws.Range($ws.Cells.Item(24, 9), $ws.Cells.Item(37, 9))
$iFillRange.Formula = "=D7/D24"

# --- Column J: constant Sa reference value (299.5941 * 0.81) repeated for
#     rows 23-37. ---
$ws.Cells.Item(23, 10).Formula = "=299.5941*0.81"
$jFillRange = $ws.Range($ws.Cells.Item(24, 10), $ws.Cells.Item(37, 10))
$jFillRange.Formula = "=299.5941*0.81"

# --- Finally select J36 to match the author's last active cell. ---
$ws.Range("J36").Select()
